$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C4"  = [double]"2.047200233e+45"
    "C5"  = 109862.141
    "C6"  = 21973.935
    "C7"  = 440.63019
    "C8"  = 76.19047399999999
    "C9"  = 2486.1433
    "C10" = 369.07993
    "C11" = 10776.1331
    "C12" = 4096.3327
    "C13" = 6694.9582
    "C16" = 137912.251
    "C17" = 498613306.088
    "C18" = 2187.9945
    "C19" = 919.12046
    "C20" = 2280712.096
    "C22" = 858.95321
    "C23" = 594.93834
    "C24" = 4640.792
    "C25" = 1277.3193
    "C26" = 1140.76412
    "C27" = 6589.9124
    "C28" = 8125.2022
    "C29" = 1550.277
    "C30" = 5035.7233
    "C31" = 2593.4862
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
